$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.NumberFormat = "General"
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "41.742.21"
$ws.Range("E2").Value = "  -1.77%  "

Set-TextValue $ws.Range("D3") "2.221.19"
$ws.Range("E3").Value = "  -1.34%  "

$ws.Range("E4").Value = "  -0.09%  "

Set-TextValue $ws.Range("D5") "249.86"
$ws.Range("E5").Value = "  +5.61%  "

$ws.Range("E6").Value = "  +0.44%  "

Set-TextValue $ws.Range("D7") "71.61"
$ws.Range("E7").Value = "  +2.16%  "

$ws.Range("E8").Value = "  -0.08%  "

Set-TextValue $ws.Range("D9") "0.609"
$ws.Range("E9").Value = "  +8.43%  "

Set-TextValue $ws.Range("D10") "40.56"
$ws.Range("E10").Value = "  +10.13%  "

$ws.Range("E11").Value = "  -3.50%  "

Set-TextValue $ws.Range("D12") "58.21"
$ws.Range("E12").Value = "  -1.84%  "

Set-TextValue $ws.Range("D13") "7.23"
$ws.Range("E13").Value = "  +6.66%  "

Set-TextValue $ws.Range("D14") "0.105"
$ws.Range("E14").Value = "  -0.84%  "

Set-TextValue $ws.Range("D15") "2.552.24"
$ws.Range("E15").Value = "  -1.13%  "

Set-TextValue $ws.Range("D16") "14.95"
$ws.Range("E16").Value = "  -1.67%  "

$ws.Range("E17").Value = "  +0.43%  "

Set-TextValue $ws.Range("D18") "2.219.09"
$ws.Range("E18").Value = "  -1.31%  "

Set-TextValue $ws.Range("D19") "41.641.09"
$ws.Range("E19").Value = "  -1.70%  "

Set-TextValue $ws.Range("D20") "0.0₃0965"
$ws.Range("E20").Value = "  -1.70%  "

Set-TextValue $ws.Range("D21") "6.22"
$ws.Range("E21").Value = "  -0.99%  "

Set-TextValue $ws.Range("D22") "72.87"
$ws.Range("E22").Value = "  -0.92%  "

Set-TextValue $ws.Range("D23") "234.42"
$ws.Range("E23").Value = "  -1.22%  "

Set-TextValue $ws.Range("D24") "2.08"
$ws.Range("E24").Value = "  +5.55%  "

$ws.Range("E25").Value = "  +9.27%  "

$ws.Range("E26").Value = "  -0.16%  "

Set-TextValue $ws.Range("D27") "2.53"
$ws.Range("E27").Value = "  +4.73%  "

$ws.Range("E28").Value = "  +7.34%  "

Set-TextValue $ws.Range("D29") "170.97"
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").Value = "  -6.44%  "

$ws.Range("E31").Value = "  +0.54%  "

Set-TextValue $ws.Range("D32") "0.123"
$ws.Range("E32").Value = "  +0.33%  "

Set-TextValue $ws.Range("D33") "5.56"
$ws.Range("E33").Value = "  +3.97%  "

$ws.Range("E34").Value = "  -2.34%  "

Set-TextValue $ws.Range("D35") "0.0733"
$ws.Range("E35").Value = "  +1.31%  "

$ws.Range("E36").Value = "  +0.08%  "

Set-TextValue $ws.Range("D37") "26.06"
$ws.Range("E37").Value = "  +16.04%  "

Set-TextValue $ws.Range("D38") "4.03"
$ws.Range("E38").Value = "  +8.92%  "

Set-TextValue $ws.Range("D39") "0.0302"
$ws.Range("E39").Value = "  +9.26%  "

$ws.Range("E40").Value = "  -0.22%  "

Set-TextValue $ws.Range("D41") "5.96"
$ws.Range("E41").Value = "  +0.02%  "

Set-TextValue $ws.Range("D42") "66.34"
$ws.Range("E42").Value = "  +1.71%  "

Set-TextValue $ws.Range("D43") "12.27"
$ws.Range("E43").Value = "  +19.41%  "

$ws.Range("E44").Value = "  +5.33%  "

Set-TextValue $ws.Range("D45") "4.85"
$ws.Range("E45").Value = "  -2.92%  "

Set-TextValue $ws.Range("D48") "4.70"
$ws.Range("E48").Value = "  +1.00%  "

$ws.Range("E49").Value = "  +0.09%  "

$ws.Range("E50").Value = "  +5.17%  "

$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D46") "0.103"
$ws.Range("E46").Value = "  -0.75%  "

$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D47") "8.69"
$ws.Range("E47").Value = "  -6.95%  "

$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "2.37"
$ws.Range("E51").Value = "  +0.54%  "
